$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-02-04 Sunday" "2024-02-05 Monday"
Replace-Text "859×2=" "845×4="
Replace-Text "661×7=" "296×9="
Replace-Text "306×7=" "435×4="
Replace-Text "489×4=" "143×5="
Replace-Text "781×3=" "918×2="
Replace-Text "502×2=" "447×2="
Replace-Text "587×4=" "684×4="
Replace-Text "359×3=" "284×3="
Replace-Text "289×4=" "614×8="
Replace-Text "195×5=" "247×2="
Replace-Text "441×3=" "524×7="
Replace-Text "399×9=" "789×3="
Replace-Text "795×9=" "413×5="
Replace-Text "991×7=" "467×6="
Replace-Text "250×7=" "548×5="
Replace-Text "966×4=" "486×6="
Replace-Text "519×4=" "310×6="
Replace-Text "163×5=" "961×3="
Replace-Text "538×8=" "270×6="
Replace-Text "819×7=" "298×2="
Replace-Text "230×6=" "904×2="
Replace-Text "133×2=" "159×7="
Replace-Text "376×5=" "543×9="
Replace-Text "493×8=" "436×8="
Replace-Text "665×7=" "741×9="
